# Apply edits described by the commit diff:
#  - Insert a new "2019" column into the BPHC sheet (before the existing
#    2020 column), shifting all later year columns one to the right.
#  - The new 2019 column's capacity cell (B2) is a formula referencing the
#    following year's value ("=C2").
#  - The sheet now extends one column further (through AH) with the last
#    data column (AG, year 2050) populated and a new trailing helper
#    column (AH) left blank but styled like the rest of the row.
#  - Give the new column a narrower custom width, like the other data
#    columns.
#  - Make BPHC the active/selected sheet (instead of "About"), with B2
#    selected.

$wb = $excel.ActiveWorkbook

$bphc = $wb.Worksheets.Item("BPHC")

# Insert a new column before column B (pushes 2020..2050 -> C..AG,
# and the formerly-empty trailing style-only cell moves to AH).
$bphc.Columns("B").Insert()

# New column B: year 2019 header, and a capacity value that mirrors
# column C (2020) via a formula.
$bphc.Range("B1").Value = 2019
$bphc.Range("B2").Formula = "=C2"

# The last data column (AG, year 2050) now carries the same capacity
# value as the rest of the row.
$bphc.Range("AG2").Value = 174

# Match the narrower column width used for the new year column
# (target stored width is 6.42578125 characters).
$bphc.Columns("B").ColumnWidth = 5.6

# Switch the active sheet to BPHC and select B2, and drop the previous
# tab selection on About.
$bphc.Activate()
$bphc.Range("B2").Select() | Out-Null
